$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O219").Value = 4417.2
$ws.Range("D220").Value = 0.00608116293155669
$ws.Range("E220").Value = 0.00523456712027648
$ws.Range("F220").Value = 0.00792756208803924
$ws.Range("G220").Value = 0.00575477607265906
$ws.Range("H220").Value = 0.0037233826809977
$ws.Range("J220").Value = 29354.3
$ws.Range("K220").Value = 1893.4
$ws.Range("N220").Value = 3139.2
$ws.Range("O220").Value = 4465.4
$ws.Range("P220").Value = 2483.3
$ws.Range("Q220").Value = 483.1
$ws.Range("S220").Value = 164.8
$ws.Range("T220").Value = 2152.299
$ws.Range("U220").Value = 221.6
$ws.Range("AG220").Value = 278.163
$ws.Range("H221").Value = 0.00545171502706143
$ws.Range("J221").Value = 29645.7262081059
$ws.Range("K221").Value = 1885.09320921534
$ws.Range("N221").Value = 3167.88900514805
$ws.Range("O221").Value = 4505.15020255048
$ws.Range("P221").Value = 2522.33520327819
$ws.Range("Q221").Value = 492.801912246684
$ws.Range("S221").Value = 163.15466849793
$ws.Range("T221").Value = 2160.527
$ws.Range("AF221").Value = 1763.0346223627
$ws.Range("AG221").Value = 280.671045427128
$ws.Range("J222").Value = 29949.8187152429
$ws.Range("K222").Value = 1907.3318350055
$ws.Range("N222").Value = 3195.56015057414
$ws.Range("O222").Value = 4578.04931402487
$ws.Range("P222").Value = 2548.44772380158
$ws.Range("Q222").Value = 502.698664280664
$ws.Range("S222").Value = 164.648939314487
$ws.Range("T222").Value = 2202.093889974
$ws.Range("AF222").Value = 1778.08535540041
$ws.Range("AG222").Value = 281.06477032852
$ws.Range("J223").Value = 30261.308341014
$ws.Range("K223").Value = 1930.4539814257
$ws.Range("N223").Value = 3223.36403051996
$ws.Range("O223").Value = 4609.71515475437
$ws.Range("P223").Value = 2574.48894008494
$ws.Range("Q223").Value = 512.794168994753
$ws.Range("S223").Value = 166.13773661157
$ws.Range("T223").Value = 2208.793889974
$ws.Range("AF223").Value = 1794.71149812467
$ws.Range("AG223").Value = 282.064144348848
$ws.Range("J224").Value = 30566.7181432502
$ws.Range("K224").Value = 1954.158923269
$ws.Range("N224").Value = 3250.16069012155
$ws.Range("O224").Value = 4641.62782503496
$ws.Range("P224").Value = 2599.52236547798
$ws.Range("Q224").Value = 523.092417862893
$ws.Range("S224").Value = 167.08465548067
$ws.Range("T224").Value = 2215.493889974
$ws.Range("AF224").Value = 1811.522377287
$ws.Range("AG224").Value = 283.064920289481
$ws.Range("J225").Value = 30856.5230650801
$ws.Range("K225").Value = 1968.11289188136
$ws.Range("N225").Value = 3276.93438119768
$ws.Range("O225").Value = 4673.78951251855
$ws.Range("P225").Value = 2624.24639104833
$ws.Range("Q225").Value = 521.771537588453
$ws.Range("S225").Value = 167.166758272789
$ws.Range("T225").Value = 2205.954889974
$ws.Range("AF225").Value = 1834.11720688667
$ws.Range("AG225").Value = 284.953140177262
$ws.Range("J226").Value = 31126.264569245
$ws.Range("K226").Value = 1980.71518106142
$ws.Range("N226").Value = 3304.82349595586
$ws.Range("O226").Value = 4755.93289067343
$ws.Range("P226").Value = 2646.38928336033
$ws.Range("Q226").Value = 520.453992718313
$ws.Range("S226").Value = 166.433306663197
$ws.Range("T226").Value = 2249.7275234174
$ws.Range("AF226").Value = 1857.01889003517
$ws.Range("AG226").Value = 286.849919297766
$ws.Range("J227").Value = 31399.8566283152
$ws.Range("K227").Value = 1995.67745272291
$ws.Range("N227").Value = 3332.32958624848
$ws.Range("O227").Value = 4788.84607274818
$ws.Range("P227").Value = 2669.39124256
$ws.Range("Q227").Value = 519.139774830117
$ws.Range("S227").Value = 166.312889234756
$ws.Range("T227").Value = 2255.6995234174
$ws.Range("AF227").Value = 1880.23203537275
$ws.Range("AG227").Value = 288.755296449687
$ws.Range("J228").Value = 31684.696360925
$ws.Range("K228").Value = 2009.40022998616
$ws.Range("N228").Value = 3360.14586941614
$ws.Range("O228").Value = 4822.01619301344
$ws.Range("P228").Value = 2692.90294087333
$ws.Range("Q228").Value = 517.828875522779
$ws.Range("S228").Value = 166.756244312196
$ws.Range("T228").Value = 2263.3995234174
$ws.Range("AF228").Value = 1903.76132307562
$ws.Range("AG228").Value = 290.669310607591
$ws.Range("J229").Value = 31972.5760053022
$ws.Range("J230").Value = 32256.3011035973
$ws.Range("J231").Value = 32537.3916116939
$ws.Range("J232").Value = 32821.8260227347
$ws.Range("J233").Value = 33110.0096582887
$ws.Range("J234").Value = 33402.5505007094
$ws.Range("J235").Value = 33702.0831401952
$ws.Range("J236").Value = 34008.1009247848
$ws.Range("J237").Value = 34320.5025240862
$ws.Range("J238").Value = 34638.3759645689
$ws.Range("J239").Value = 34961.8225766253
$ws.Range("J240").Value = 35290.0317171173
$ws.Range("J241").Value = 35622.0914125148
$ws.Range("J242").Value = 35957.1910196797
$ws.Range("J243").Value = 36295.8371905733
$ws.Range("J244").Value = 36638.1312555879
$ws.Range("J245").Value = 36983.9718843311
$ws.Range("J246").Value = 37332.953755234
$ws.Range("J247").Value = 37685.1781986889
$ws.Range("J248").Value = 38040.8478754802
$ws.Range("J249").Value = 38398.9494816855
$ws.Range("J250").Value = 38759.4830173047
$ws.Range("J251").Value = 39123.0564646914
$ws.Range("J252").Value = 39489.87248463
$ws.Range("J253").Value = 39859.5257555516
$ws.Range("J254").Value = 40231.9149470639
$ws.Range("J255").Value = 40607.2427199514
$ws.Range("J256").Value = 40985.7117349985
$ws.Range("J257").Value = 41367.4233225977
$ws.Range("J258").Value = 41752.6814739255
$ws.Range("J259").Value = 42141.2835281975
$ws.Range("J260").Value = 42533.4321461982
